# Generate Report for Handoff
#
# Marks the file "6b07ccfe-699f-489d-aef5-2635bc7b185d.md" as ready for
# handoff in both locale sheets (zh-cn, de-de) and on the Overview sheet,
# and stamps the new "Latest Handoff Datetime" for that handoff.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: update the Status column for the target file row ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: update Status + Latest Handoff Datetime for the target file row ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-02-22 13:53:11"

# --- de-de sheet: update Status + Latest Handoff Datetime for the target file row ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-02-22 13:53:26"
